# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 172 (pushing the existing
# rows 172:190 down to 173:191) on the single data sheet, then populate
# the new row with the latest report's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 172:190 down to 173:191, leaving a blank row 172 behind
# (row formatting/style is inherited from the surrounding rows, same as
# Excel's native Insert behavior).
$ws.Rows("172:172").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A172").Value = 11
$ws.Range("B172").Value = "Vega Monumental Concepción"
$ws.Range("C172").Value = "Bíobío"
$ws.Range("D172").Value = 44776
$ws.Range("E172").Value = 8
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108005
$ws.Range("J172").Value = "Piña"
$ws.Range("K172").Value = "Caramelo"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 200
$ws.Range("N172").Value = 19000
$ws.Range("O172").Value = 20000
$ws.Range("P172").Value = 19500
$ws.Range("Q172").Value = "$/caja 12 unidades"
$ws.Range("R172").Value = "Ecuador"
$ws.Range("S172").Value = 1625
$ws.Range("T172").Value = 12
